$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B entirely; this shifts column C (quantities) into B,
# matching the target layout where the old "3 sạch 17_09_2016" / quantity
# column is gone and the HQ Food / quantity values that used to live in C
# now live in B.
$ws.Range("B:B").Delete()

# Update the timestamp cell with the new value.
$ws.Range("B1").Value = "Wed Sep 21 2016 15:31:21 GMT+0700 (SE Asia Standard Time)"
